# #5: property boat&car done
# Extend the "汽車" (car) sheet from 7 columns (A:G) to 14 columns (A:N),
# matching the header/meta layout already used by the other property
# sheets (name/capacity/owner/register_date/register_reason/
# acquire_value/property_category/category/date/legislator_name/
# legislator_id/source_file/index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")
$land = $wb.Worksheets.Item("土地")

# --- Row 1: header labels (B1:N1) ---------------------------------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Columns H1:N1 need the same bold/bordered header formatting as the rest
# of row 1 -- copy formats only from an existing header cell.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: data (H2:N2) -------------------------------------------------
# These columns mirror the meta columns already present on the other
# property sheets (property_category/category/date/legislator_name/
# legislator_id/source_file/index) -- copy the values straight from the
# "土地" (land) sheet's first data row, which carries identical metadata
# for this same filing, then fix up the row-specific index (N2).
$land.Range("I2:N2").Copy()
$ws.Range("H2").PasteSpecial(-4104)
$excel.CutCopyMode = $false
$ws.Cells.Item(2, 14).Value = 29
